$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
# B1 currently holds "教育得分"; split it into separate indicator columns.
$ws.Range("B1").Value = "教育"
$ws.Range("C1").Value = "交通"
$ws.Range("D1").Value = "医疗"
$ws.Range("E1").Value = "得分"

# Give the new header cells (C1, D1, E1) the same formatting as B1
# (bold font, thin border, centered) by copying B1's format.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- City data: Name, 教育(B, already present), 交通(C), 医疗(D) ---
# 得分(E) = average of 教育/交通/医疗
$data = @(
    @("上海市", 6.8, 5.2, 5.2),
    @("合肥市", 5, 3.4, 3.4),
    @("安庆市", 2.3, 3.8, 3.8),
    @("宣城市", 2.3, 3.6, 3.6),
    @("池州市", 4.2, 6.2, 6.2),
    @("滁州市", 3.5, 4.4, 4.4),
    @("芜湖市", 4, 4.3, 4.3),
    @("铜陵市", 2.1, 1, 1),
    @("马鞍山市", 2, 1.1, 1.1),
    @("南京市", 5.3, 4.2, 4.2),
    @("南通市", 3.4, 3.5, 3.5),
    @("常州市", 2.4, 4.1, 4.1),
    @("扬州市", 2.8, 2.4, 2.4),
    @("无锡市", 2, 3.6, 3.6),
    @("泰州市", 2.3, 4, 4),
    @("盐城市", 3.1, 3, 3),
    @("苏州市", 2.7, 2.6, 2.6),
    @("镇江市", 3.5, 3.3, 3.3),
    @("台州市", 1.6, 3.5, 3.5),
    @("嘉兴市", 3.2, 1.3, 1.3),
    @("宁波市", 3.2, 2.2, 2.2),
    @("杭州市", 4.9, 1.6, 1.6),
    @("温州市", 3, 1.8, 1.8),
    @("湖州市", 1.6, 3.7, 3.7),
    @("绍兴市", 3.3, 0.9, 0.9),
    @("舟山市", 3.3, 5.4, 5.4),
    @("金华市", 4.4, 4.1, 4.1)
)

$row = 2
foreach ($item in $data) {
    $b = $item[1]
    $c = $item[2]
    $d = $item[3]

    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = ($b + $c + $d) / 3

    $row = $row + 1
}
